$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.871715666666666
$ws.Range("H2").Value = 23.615147
$ws.Range("I2").Value = 0.02771913691218268
$ws.Range("J2").Value = 0.02771913691218268
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 268.1952276871258
$ws.Range("R2").Value = 2413.757049184132
$ws.Range("S2").Value = 0.01433485323773431
$ws.Range("T2").Value = 0.01433485323773431

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.871715666666666
$ws.Range("H3").Value = 23.615147
$ws.Range("I3").Value = 0.02771913691218268
$ws.Range("J3").Value = 0.02771913691218268
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.05649099999999
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 217.932360474353
$ws.Range("R3").Value = 1961.391244269177
$ws.Range("S3").Value = 0.01164833703453263
$ws.Range("T3").Value = 0.01164833703453263

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.871715666666666
$ws.Range("H4").Value = 23.615147
$ws.Range("I4").Value = 0.02771913691218268
$ws.Range("J4").Value = 0.02771913691218268
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 32.47836560470344
$ws.Range("R4").Value = 292.3052904423309
$ws.Range("S4").Value = 0.001735946639915736
$ws.Range("T4").Value = 0.001735946639915736

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 266.1315866666666
$ws.Range("H5").Value = 798.3947599999999
$ws.Range("I5").Value = 0.9371448614065047
$ws.Range("J5").Value = 0.9371448614065045
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 9067.301780607513
$ws.Range("R5").Value = 81605.71602546761
$ws.Range("S5").Value = 0.4846411377568857
$ws.Range("T5").Value = 0.4846411377568858

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 266.1315866666666
$ws.Range("H6").Value = 798.3947599999999
$ws.Range("I6").Value = 0.9371448614065047
$ws.Range("J6").Value = 0.9371448614065045
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.05649099999999
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 7367.985244265238
$ws.Range("R6").Value = 66311.86719838715
$ws.Range("S6").Value = 0.3938138200488354
$ws.Range("T6").Value = 0.3938138200488354

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 266.1315866666666
$ws.Range("H7").Value = 798.3947599999999
$ws.Range("I7").Value = 0.9371448614065047
$ws.Range("J7").Value = 0.9371448614065045
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 1098.047660349498
$ws.Range("R7").Value = 9882.428943145478
$ws.Range("S7").Value = 0.05868990360078345
$ws.Range("T7").Value = 0.05868990360078345

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.977966333333333
$ws.Range("H8").Value = 29.933899
$ws.Range("I8").Value = 0.03513600168131278
$ws.Range("J8").Value = 0.03513600168131277
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 339.95675986554
$ws.Range("R8").Value = 3059.61083878986
$ws.Range("S8").Value = 0.01817045851961717
$ws.Range("T8").Value = 0.01817045851961717

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.977966333333333
$ws.Range("H9").Value = 29.933899
$ws.Range("I9").Value = 0.03513600168131278
$ws.Range("J9").Value = 0.03513600168131277
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.05649099999999
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 276.244956987601
$ws.Range("R9").Value = 2486.204612888409
$ws.Range("S9").Value = 0.01476510581575712
$ws.Range("T9").Value = 0.01476510581575712

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.977966333333333
$ws.Range("H10").Value = 29.933899
$ws.Range("I10").Value = 0.03513600168131278
$ws.Range("J10").Value = 0.03513600168131277
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 41.16866669075855
$ws.Range("R10").Value = 370.5180002168269
$ws.Range("S10").Value = 0.002200437345938478
$ws.Range("T10").Value = 0.002200437345938478
